$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as exact text, avoiding Excel auto-numeric-conversion
# for values that look like numbers (keeps trailing zeros / exact formatting).
function Set-TextValue($range, [string]$text) {
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $range.NumberFormat = "@"
    }
    $range.Value = $text
}

# --- Update price (D) and volume (E) values for rows with unchanged coin identity ---
$ws.Range("D2").Value = "41.870.02"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.219.95"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "251.06"
$ws.Range("E5").Value = "  -1.27%  "
Set-TextValue $ws.Range("D6") "0.624"
$ws.Range("E6").Value = "  -0.51%  "
Set-TextValue $ws.Range("D7") "68.12"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.634"
$ws.Range("E9").Value = "  +7.11%  "
Set-TextValue $ws.Range("D10") "39.64"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("E11").Value = "  +2.90%  "
Set-TextValue $ws.Range("D12") "0.0937"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "2.544.85"
$ws.Range("E15").Value = "  +0.93%  "
Set-TextValue $ws.Range("D16") "14.68"
$ws.Range("E16").Value = "  -0.88%  "
Set-TextValue $ws.Range("D17") "0.875"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "2.206.63"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "41.790.82"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("E20").Value = "  +0.37%  "
Set-TextValue $ws.Range("D21") "6.24"
$ws.Range("E21").Value = "  -0.35%  "
Set-TextValue $ws.Range("D22") "72.75"
$ws.Range("E22").Value = "  +1.12%  "
Set-TextValue $ws.Range("D23") "232.99"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("E28").Value = "  -4.22%  "
Set-TextValue $ws.Range("D29") "3.70"
$ws.Range("E29").Value = "  -1.75%  "
Set-TextValue $ws.Range("D30") "2.15"
$ws.Range("E30").Value = "  -1.76%  "
Set-TextValue $ws.Range("D31") "167.29"
$ws.Range("E31").Value = "  -1.77%  "
Set-TextValue $ws.Range("D32") "20.42"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E35").Value = "  -0.51%  "
Set-TextValue $ws.Range("D36") "0.123"
$ws.Range("E36").Value = "  -0.53%  "
Set-TextValue $ws.Range("D37") "4.63"
$ws.Range("E37").Value = "  -0.24%  "
Set-TextValue $ws.Range("D38") "4.13"
$ws.Range("E38").Value = "  +3.04%  "
Set-TextValue $ws.Range("D39") "25.58"
$ws.Range("E39").Value = "  -2.05%  "
Set-TextValue $ws.Range("D40") "0.0308"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("E41").Value = "  +0.45%  "
Set-TextValue $ws.Range("D42") "12.12"
$ws.Range("E42").Value = "  +0.34%  "
Set-TextValue $ws.Range("D43") "5.65"
$ws.Range("E43").Value = "  -2.22%  "
Set-TextValue $ws.Range("D44") "5.15"
$ws.Range("E44").Value = "  +3.06%  "
Set-TextValue $ws.Range("D45") "62.24"
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("E46").Value = "  -1.33%  "
Set-TextValue $ws.Range("D47") "8.59"
$ws.Range("E47").Value = "  -1.14%  "
Set-TextValue $ws.Range("D48") "0.0999"
$ws.Range("E48").Value = "  -0.81%  "
Set-TextValue $ws.Range("D49") "0.999"
$ws.Range("E49").Value = "  -0.50%  "
Set-TextValue $ws.Range("D50") "1.17"
$ws.Range("E50").Value = "  +0.72%  "

# --- Rows where the coin identity itself changed (ranking reshuffled) ---
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "11.37"
$ws.Range("E26").Value = "  -5.14%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D33") "5.96"
$ws.Range("E33").Value = "  +7.23%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.0798"
$ws.Range("E34").Value = "  +8.92%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws.Range("D51") "4.34"
$ws.Range("E51").Value = "  +0.42%  "

